$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the BMSMake/Udaan column (I1:I2), shifting RGP/IsActive (J,K) left to I,J.
# This also correctly updates the dataValidation sqref and sheet dimension.
$ws.Range("I1:I2").Delete(-4159)

# At this point the layout is:
# A=SpareCategory B=PartCode C=Description D=UOM E=MinQty F=AvailableQty G=TentativeCost H=ProductMake I=RGP J=IsActive
# Target layout is:
# A=SpareCategory B=ProductMake C=SparePartCode D=SparePartDescription E=UOM F=MinQty G=AvailableQty H=TentativeCost I=RGP J=IsActive
# Rewrite B:H directly (content-only, no structural column moves) to avoid disturbing <cols>/bestFit metadata.

$ws.Range("B1").Value = "ProductMake"
$ws.Range("B2").Value = "FBTECH"

$ws.Range("C1").Value = "SparePartCode"
$ws.Range("C2").Value = 123

$ws.Range("D1").Value = "SparePartDescription"
$ws.Range("D2").Value = "Hard Disk"

$ws.Range("E1").Value = "UOM"
$ws.Range("E2").Value = "KG"

$ws.Range("F1").Value = "MinQty"
$ws.Range("F2").Value = 1

$ws.Range("G1").Value = "AvailableQty"
$ws.Range("G2").Value = 10

$ws.Range("H1").Value = "TentativeCost"
$ws.Range("H2").Value = 10

# Column widths: column A (SpareCategory) widened; ProductMake's former width/bestFit
# (now in column B) and the rest keep their pre-existing <col> definitions.
$ws.Range("A1").EntireColumn.ColumnWidth = 13.166666666666666

# Fix selection to match the target workbook state.
$ws.Range("I19").Select()
